$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "315.73"
$ws.Range("D2").Style = $origStyle

$origStyle = $ws.Range("E2").Style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.33%"
$ws.Range("E2").Style = $origStyle

$origStyle = $ws.Range("E3").Style
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.04%"
$ws.Range("E3").Style = $origStyle

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.165"
$ws.Range("D4").Style = $origStyle

$origStyle = $ws.Range("E4").Style
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.63%"
$ws.Range("E4").Style = $origStyle

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08095"
$ws.Range("D5").Style = $origStyle

$origStyle = $ws.Range("E5").Style
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.89%"
$ws.Range("E5").Style = $origStyle

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.523"
$ws.Range("D6").Style = $origStyle

$origStyle = $ws.Range("E6").Style
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.95%"
$ws.Range("E6").Style = $origStyle

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.682"
$ws.Range("D7").Style = $origStyle

$origStyle = $ws.Range("E7").Style
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.57%"
$ws.Range("E7").Style = $origStyle

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.090"
$ws.Range("D8").Style = $origStyle

$origStyle = $ws.Range("E8").Style
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "17.02%"
$ws.Range("E8").Style = $origStyle

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1304"
$ws.Range("D9").Style = $origStyle

$origStyle = $ws.Range("E9").Style
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8.51%"
$ws.Range("E9").Style = $origStyle

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1932"
$ws.Range("D10").Style = $origStyle

$origStyle = $ws.Range("E10").Style
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.29%"
$ws.Range("E10").Style = $origStyle

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09449"
$ws.Range("D11").Style = $origStyle

$origStyle = $ws.Range("E11").Style
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.51%"
$ws.Range("E11").Style = $origStyle

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04272"
$ws.Range("D12").Style = $origStyle

$origStyle = $ws.Range("E12").Style
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "7.89%"
$ws.Range("E12").Style = $origStyle

$origStyle = $ws.Range("E13").Style
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.79%"
$ws.Range("E13").Style = $origStyle

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001315"
$ws.Range("D14").Style = $origStyle

$origStyle = $ws.Range("E14").Style
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.45%"
$ws.Range("E14").Style = $origStyle

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005869"
$ws.Range("D15").Style = $origStyle

$origStyle = $ws.Range("E15").Style
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.34%"
$ws.Range("E15").Style = $origStyle

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.399"
$ws.Range("D17").Style = $origStyle

$origStyle = $ws.Range("E17").Style
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.02%"
$ws.Range("E17").Style = $origStyle

$origStyle = $ws.Range("E18").Style
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.09%"
$ws.Range("E18").Style = $origStyle

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3367"
$ws.Range("D19").Style = $origStyle

$origStyle = $ws.Range("E19").Style
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.42%"
$ws.Range("E19").Style = $origStyle

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.310"
$ws.Range("D20").Style = $origStyle

$origStyle = $ws.Range("E20").Style
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.80%"
$ws.Range("E20").Style = $origStyle

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1385"
$ws.Range("D21").Style = $origStyle

$origStyle = $ws.Range("E21").Style
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.33%"
$ws.Range("E21").Style = $origStyle

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3146"
$ws.Range("D22").Style = $origStyle

$origStyle = $ws.Range("E22").Style
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.89%"
$ws.Range("E22").Style = $origStyle

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04258"
$ws.Range("D23").Style = $origStyle

$origStyle = $ws.Range("E23").Style
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.73%"
$ws.Range("E23").Style = $origStyle

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001280"
$ws.Range("D24").Style = $origStyle

$origStyle = $ws.Range("E24").Style
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.27%"
$ws.Range("E24").Style = $origStyle

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004247"
$ws.Range("D25").Style = $origStyle

$origStyle = $ws.Range("E25").Style
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.57%"
$ws.Range("E25").Style = $origStyle

$origStyle = $ws.Range("E26").Style
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.34%"
$ws.Range("E26").Style = $origStyle

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02708"
$ws.Range("D38").Style = $origStyle

$origStyle = $ws.Range("E38").Style
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "12.00%"
$ws.Range("E38").Style = $origStyle

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05465"
$ws.Range("D39").Style = $origStyle

$origStyle = $ws.Range("E39").Style
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.85%"
$ws.Range("E39").Style = $origStyle

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005440"
$ws.Range("D40").Style = $origStyle

$origStyle = $ws.Range("E40").Style
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-9.51%"
$ws.Range("E40").Style = $origStyle

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007771"
$ws.Range("D41").Style = $origStyle

$origStyle = $ws.Range("E41").Style
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.16%"
$ws.Range("E41").Style = $origStyle

$origStyle = $ws.Range("E42").Style
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.82%"
$ws.Range("E42").Style = $origStyle

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007373"
$ws.Range("D43").Style = $origStyle

$origStyle = $ws.Range("E43").Style
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.16%"
$ws.Range("E43").Style = $origStyle

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008574"
$ws.Range("D44").Style = $origStyle

$origStyle = $ws.Range("E44").Style
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "18.67%"
$ws.Range("E44").Style = $origStyle

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3144"
$ws.Range("D45").Style = $origStyle

$origStyle = $ws.Range("E45").Style
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.52%"
$ws.Range("E45").Style = $origStyle

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006797"
$ws.Range("D46").Style = $origStyle

$origStyle = $ws.Range("E47").Style
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.40%"
$ws.Range("E47").Style = $origStyle

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06257"
$ws.Range("D48").Style = $origStyle

$origStyle = $ws.Range("E48").Style
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "36.61%"
$ws.Range("E48").Style = $origStyle

$origStyle = $ws.Range("E49").Style
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.11%"
$ws.Range("E49").Style = $origStyle

$origStyle = $ws.Range("E50").Style
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.40%"
$ws.Range("E50").Style = $origStyle

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001992"
$ws.Range("D51").Style = $origStyle

$origStyle = $ws.Range("E51").Style
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.40%"
$ws.Range("E51").Style = $origStyle
